$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: IonQ, Inc. / IONQ
$ws.Range("D2").Value = 54.76
$ws.Range("E2").Value = 65.3
$ws.Range("F2").Value = 16.76
$ws.Range("K2").Value = 58.5
$ws.Range("N2").Value = 54.85170003294819

# Row 3: swaps from Rigetti/RGTI -> D-Wave Quantum Inc./QBTS
$ws.Range("B3").Value = "D-Wave Quantum Inc."
$ws.Range("C3").Value = "QBTS"
$ws.Range("D3").Value = 28.73
$ws.Range("E3").Value = 65.8
$ws.Range("F3").Value = 28.2
$ws.Range("H3").Value = 70
$ws.Range("J3").Value = 76
$ws.Range("K3").Value = 56.7
$ws.Range("N3").Value = 54.85170003294819

# Row 4: swaps from D-Wave Quantum Inc./QBTS -> Rigetti Computing, Inc./RGTI
$ws.Range("B4").Value = "Rigetti Computing, Inc."
$ws.Range("C4").Value = "RGTI"
$ws.Range("D4").Value = 30.06
$ws.Range("E4").Value = 63.1
$ws.Range("F4").Value = 17.56
$ws.Range("H4").Value = 63
$ws.Range("I4").Value = 60
$ws.Range("J4").Value = 86
$ws.Range("K4").Value = 55.5
$ws.Range("N4").Value = 54.85170003294819

# Row 5: International Business Machines / IBM
$ws.Range("D5").Value = 307.99
$ws.Range("E5").Value = 52.9
$ws.Range("F5").Value = 1.58
$ws.Range("I5").Value = 60
$ws.Range("J5").Value = 56
$ws.Range("K5").Value = 55.5
$ws.Range("N5").Value = 54.85170003294819
